$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update likely_death_date values for individuals F59892 (row 3) and F59893 (row 4)
$ws.Range("B3").Value = 44830
$ws.Range("B4").Value = 44834
